# Auto commit at 2025-08-28 10:28:07.77
#
# Update the "Metrics" sheet's monthly/yearly/total figures (B2:B13) with the
# latest numbers, which ripples into the "today" sheet via its formulas
# (Metrics!B2..B13), and move the active selection around to match where the
# author ended up working.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metrics sheet: refresh the reported values in column B.
# ---------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsMetrics.Activate() | Out-Null

$wsMetrics.Range("B2").Value  = 456313.11
$wsMetrics.Range("B3").Value  = 391591.01
$wsMetrics.Range("B4").Value  = 144298.21
$wsMetrics.Range("B5").Value  = 17898
$wsMetrics.Range("B6").Value  = 3851941.6799999997
$wsMetrics.Range("B7").Value  = 3270305.67
$wsMetrics.Range("B8").Value  = 1107940.77
$wsMetrics.Range("B9").Value  = 148586
$wsMetrics.Range("B10").Value = 32317265.480999827
$wsMetrics.Range("B11").Value = 19300175.740000002
$wsMetrics.Range("B12").Value = 11389649.660000002
$wsMetrics.Range("B13").Value = 1246213

# Leave the selection on Metrics where the author left it.
$wsMetrics.Range("G12").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. IncomeChart sheet: no data changes, just passed through as the active
#    tab moves elsewhere below (its tabSelected flag is cleared as a side
#    effect of another sheet becoming active).
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 3. today sheet: becomes the active tab, with a new selection. Its formulas
#    referencing Metrics!B2:B13 (and the downstream E/F columns) recalculate
#    automatically.
# ---------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Activate() | Out-Null
$wsToday.Range("F19").Select() | Out-Null
